$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Positve Testdata"
#   Row 2/3 test data replaced with new fake values; password/confirm cells
#   now both read "Test@123"; B2/B3 hyperlinks repointed to the new emails.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "fake22101"
$ws1.Range("B2").Value = "fake22101@g.com"
$ws1.Range("C2").Value = 1267432221
$ws1.Range("D2").Value = "Test@123"
$ws1.Range("E2").Value = "Test@123"

$ws1.Range("A3").Value = "fake22071"
$ws1.Range("B3").Value = "fake22071@g.com"
$ws1.Range("C3").Value = 3456789044
$ws1.Range("D3").Value = "Test@123"
$ws1.Range("E3").Value = "Test@123"

# Drop every existing hyperlink on the sheet, then recreate only the two
# e-mail links (B2, B3) against the freshly typed addresses.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "mailto:fake22101@g.com")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "mailto:fake22071@g.com")
$ws1.Range("B2").Style = "Hyperlink"
$ws1.Range("B3").Style = "Hyperlink"

[void]$ws1.Range("F6").Select()

# ---------------------------------------------------------------------------
# Sheet 2: "Negative Testdata"
#   Row 2's phone number changes, a brand-new row 4 is appended, and every
#   hyperlink on the sheet is removed.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C2").Value = 9456789021

$ws2.Range("A4").Value = "fake3"
$ws2.Range("B4").Value = "new324@g.com"
$ws2.Range("C4").Value = 123
$ws2.Range("D4").Value = 123
$ws2.Range("E4").Value = 123

$ws2.Hyperlinks.Delete()

[void]$ws2.Rows("5:5").Select()

# Leave "Positve Testdata" as the active/selected tab, matching the saved file.
[void]$ws1.Activate()
[void]$ws1.Range("F6").Select()
